# Portugal Primeira Liga - base update (12-04-2024 20:28)
#
# The underlying data rows were re-sorted/updated. Net effect on the sheet:
#  1) 17 pairs of adjacent rows swap their B:AC content (the leading index
#     column A, which is just the running row number, stays put).
#  2) One row (253) was removed from the dataset: rows 254-261 shift their
#     B:AC content up by one row, and the now-duplicate trailing row (261)
#     is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($RowA, $RowB) {
    $rangeA = $ws.Range("B$RowA`:AC$RowA")
    $rangeB = $ws.Range("B$RowB`:AC$RowB")
    $dataA = $rangeA.Value2
    $dataB = $rangeB.Value2
    $rangeA.Value2 = $dataB
    $rangeB.Value2 = $dataA
}

# 1) Swap the 17 adjacent row pairs (B:AC only - column A is untouched).
$rowPairs = @(
    @(30,31),
    @(48,49),
    @(70,71),
    @(87,88),
    @(97,98),
    @(123,124),
    @(128,129),
    @(133,134),
    @(139,140),
    @(151,152),
    @(167,168),
    @(177,178),
    @(195,196),
    @(212,213),
    @(220,221),
    @(238,239),
    @(245,246)
)

foreach ($pair in $rowPairs) {
    Swap-RowData $pair[0] $pair[1]
}

# 2) Remove row 253's data: shift B:AC of rows 254..261 up into 253..260,
#    then drop the now-empty trailing row 261 entirely.
for ($r = 253; $r -le 260; $r++) {
    $src = $ws.Range("B$($r+1):AC$($r+1)").Value2
    $ws.Range("B$r`:AC$r").Value2 = $src
}

$ws.Range("B261:AC261").ClearContents()
$ws.Rows("261:261").Delete(-4162) | Out-Null
